# Scheduled runner update: refresh market-price-derived columns (H-N)
# across the Leve profit sheets, per upstream data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 21381.092
$ws.Range("J19").Value = 38258
$ws.Range("L19").Value = 38258
$ws.Range("N19").Value = -38608
$ws.Range("H51").Value = 4000
$ws.Range("I51").Value = 4000
$ws.Range("J51").Value = 4000
$ws.Range("K51").Value = 4000
$ws.Range("L51").Value = 4000
$ws.Range("M51").Value = -3516
$ws.Range("N51").Value = -4968
$ws.Range("H70").Value = 1993.75
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 1975
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 5925
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -6465
$ws.Range("H73").Value = 1993.75
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 1975
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 5925
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -7797
$ws.Range("H129").Value = 421.22223
$ws.Range("I129").Value = 423.875
$ws.Range("J129").Value = 400
$ws.Range("K129").Value = 1271.625
$ws.Range("L129").Value = 1200
$ws.Range("M129").Value = 3728.375
$ws.Range("N129").Value = -11200
$ws.Range("H132").Value = 6107.207
$ws.Range("J132").Value = 13966.714
$ws.Range("L132").Value = 41900.142
$ws.Range("N132").Value = -46960.142
$ws.Range("H138").Value = 2791.8936
$ws.Range("I138").Value = 1908.9048
$ws.Range("J138").Value = 3505.077
$ws.Range("K138").Value = 5726.7144
$ws.Range("L138").Value = 10515.231
$ws.Range("M138").Value = -586.7143999999998
$ws.Range("N138").Value = -20795.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 469.75
$ws.Range("I16").Value = 469.75
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 469.75
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -182.75
$ws.Range("N16").ClearContents()
$ws.Range("H21").Value = 3877.3845
$ws.Range("J21").Value = 3531.1667
$ws.Range("L21").Value = 3531.1667
$ws.Range("N21").Value = -4279.1667
$ws.Range("H28").Value = 24304.285
$ws.Range("J28").Value = 29818.182
$ws.Range("L28").Value = 29818.182
$ws.Range("N28").Value = -30202.182
$ws.Range("H63").Value = 3420.7646
$ws.Range("J63").Value = 3433.697
$ws.Range("L63").Value = 3433.697
$ws.Range("N63").Value = -4805.697
$ws.Range("H66").Value = 3420.7646
$ws.Range("J66").Value = 3433.697
$ws.Range("L66").Value = 17168.485
$ws.Range("N66").Value = -24032.485
$ws.Range("H97").Value = 45501376
$ws.Range("I97").Value = 50001510
$ws.Range("K97").Value = 50001510
$ws.Range("M97").Value = -50001014
$ws.Range("H99").Value = 24304.285
$ws.Range("J99").Value = 29818.182
$ws.Range("L99").Value = 29818.182
$ws.Range("N99").Value = -35808.182
$ws.Range("H131").Value = 93000
$ws.Range("J131").Value = 93000
$ws.Range("L131").Value = 93000
$ws.Range("N131").Value = -103080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 41049.4
$ws.Range("J82").Value = 49249.25
$ws.Range("L82").Value = 49249.25
$ws.Range("N82").Value = -50015.25
$ws.Range("H85").Value = 41049.4
$ws.Range("J85").Value = 49249.25
$ws.Range("L85").Value = 49249.25
$ws.Range("N85").Value = -51901.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1832.8158
$ws.Range("I16").Value = 1840.8077
$ws.Range("J16").Value = 1815.5
$ws.Range("K16").Value = 1840.8077
$ws.Range("L16").Value = 1815.5
$ws.Range("M16").Value = -1553.8077
$ws.Range("N16").Value = -2389.5
$ws.Range("H51").Value = 74000
$ws.Range("J51").Value = 74000
$ws.Range("L51").Value = 74000
$ws.Range("N51").Value = -75472
$ws.Range("H61").Value = 74000
$ws.Range("J61").Value = 74000
$ws.Range("L61").Value = 74000
$ws.Range("N61").Value = -74696
$ws.Range("H113").Value = 1832.8158
$ws.Range("I113").Value = 1840.8077
$ws.Range("J113").Value = 1815.5
$ws.Range("K113").Value = 1840.8077
$ws.Range("L113").Value = 1815.5
$ws.Range("M113").Value = 329.1922999999999
$ws.Range("N113").Value = -6155.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 994.7568
$ws.Range("I5").Value = 767.8333
$ws.Range("J5").Value = 1967.2858
$ws.Range("K5").Value = 2303.4999
$ws.Range("L5").Value = 5901.857400000001
$ws.Range("M5").Value = -2191.4999
$ws.Range("N5").Value = -6125.857400000001
$ws.Range("H9").Value = 673702.0600000001
$ws.Range("I9").Value = 1233604.4
$ws.Range("J9").Value = 113799.78
$ws.Range("K9").Value = 3700813.2
$ws.Range("L9").Value = 341399.34
$ws.Range("M9").Value = -3700589.2
$ws.Range("N9").Value = -341847.34
$ws.Range("H33").Value = 315.6154
$ws.Range("I33").Value = 212
$ws.Range("J33").Value = 346.7
$ws.Range("K33").Value = 1272
$ws.Range("L33").Value = 2080.2
$ws.Range("M33").Value = -989
$ws.Range("N33").Value = -2646.2
$ws.Range("H47").Value = 543.8
$ws.Range("I47").Value = 179.75
$ws.Range("K47").Value = 539.25
$ws.Range("M47").Value = -108.25
$ws.Range("H131").Value = 7356198.5
$ws.Range("I131").Value = 14707397
$ws.Range("K131").Value = 44122191
$ws.Range("M131").Value = -44117151
$ws.Range("H135").Value = 994.7568
$ws.Range("I135").Value = 767.8333
$ws.Range("J135").Value = 1967.2858
$ws.Range("K135").Value = 6910.4997
$ws.Range("L135").Value = 17705.5722
$ws.Range("M135").Value = -4375.4997
$ws.Range("N135").Value = -22775.5722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 258783.58
$ws.Range("J3").Value = 9582.091
$ws.Range("L3").Value = 9582.091
$ws.Range("N3").Value = -9814.091
$ws.Range("H11").Value = 1266825.5
$ws.Range("I11").Value = 5007500
$ws.Range("J11").Value = 19934
$ws.Range("K11").Value = 5007500
$ws.Range("L11").Value = 19934
$ws.Range("M11").Value = -5007361
$ws.Range("N11").Value = -20212
$ws.Range("H12").Value = 29999.5
$ws.Range("I12").Value = 29999.5
$ws.Range("K12").Value = 29999.5
$ws.Range("M12").Value = -29859.5
$ws.Range("H33").Value = 19000
$ws.Range("J33").Value = 19000
$ws.Range("L33").Value = 19000
$ws.Range("N33").Value = -19504
$ws.Range("H102").Value = 5599.892
$ws.Range("I102").Value = 3607.3667
$ws.Range("J102").Value = 14139.286
$ws.Range("K102").Value = 3607.3667
$ws.Range("L102").Value = 14139.286
$ws.Range("M102").Value = -1985.3667
$ws.Range("N102").Value = -17383.286
$ws.Range("H132").Value = 4950.9536
$ws.Range("I132").Value = 2780.5962
$ws.Range("J132").Value = 13632.385
$ws.Range("K132").Value = 8341.7886
$ws.Range("L132").Value = 40897.155
$ws.Range("M132").Value = -5811.7886
$ws.Range("N132").Value = -45957.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7011.25
$ws.Range("I40").Value = 6012.857
$ws.Range("K40").Value = 6012.857
$ws.Range("M40").Value = -5876.857
$ws.Range("H46").Value = 2437.9
$ws.Range("I46").Value = 911
$ws.Range("J46").Value = 3092.2856
$ws.Range("K46").Value = 911
$ws.Range("L46").Value = 3092.2856
$ws.Range("M46").Value = -723
$ws.Range("N46").Value = -3468.2856
$ws.Range("H68").Value = 9527350
$ws.Range("J68").Value = 26500
$ws.Range("L68").Value = 26500
$ws.Range("N68").Value = -27998
$ws.Range("H71").Value = 9527350
$ws.Range("J71").Value = 26500
$ws.Range("L71").Value = 132500
$ws.Range("N71").Value = -139988
$ws.Range("H132").Value = 15516.489
$ws.Range("J132").Value = 2841.8572
$ws.Range("L132").Value = 8525.571599999999
$ws.Range("N132").Value = -13585.5716
$ws.Range("H136").Value = 6003462
$ws.Range("I136").Value = 8184383.5
$ws.Range("J136").Value = 5928.125
$ws.Range("K136").Value = 24553150.5
$ws.Range("L136").Value = 17784.375
$ws.Range("M136").Value = -24550600.5
$ws.Range("N136").Value = -22884.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2225.2983
$ws.Range("I132").Value = 2380.2046
$ws.Range("K132").Value = 7140.6138
$ws.Range("M132").Value = -4610.6138
$ws.Range("H137").Value = 149973.25
$ws.Range("J137").Value = 149973.25
$ws.Range("L137").Value = 149973.25
$ws.Range("N137").Value = -160173.25
